# Apply updated "dSF" (column F) values to the matching rows on Sheet1.
# The diff only touches column F values; everything else stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = 4
    6  = 3
    11 = 1
    18 = -1
    24 = 0
    32 = -4
    34 = 3
    35 = 1
    39 = 4
    40 = 0
    46 = -5
    49 = 2
    52 = 1
    60 = -2
    62 = -2
    64 = 3
    66 = 1
    68 = -1
    70 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
